# "Alteracao fPedido - QTD"
# Fills in the "QTD" (Qtde Estoque, column E) evaluation points that were
# already present for column C, and brings the sheet's view/scroll state up
# to date with where the author was working (row 8 / around E22).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avaliação")

# --- Column E ("Qtde Estoque") scores, mirroring column C's values -------
$ws.Range("E3").Value  = 1
$ws.Range("E4").Value  = 1
$ws.Range("E5").Value  = 2
$ws.Range("E6").Value  = 3
$ws.Range("E7").Value  = 2
$ws.Range("E8").Value  = 1
$ws.Range("E9").Value  = 1
$ws.Range("E20").Value = 1
$ws.Range("E21").Value = 1

# --- Row 8 grew a touch taller (matches the other wrapped-text rows) -----
$ws.Rows.Item(8).RowHeight = 25.5

# --- Scroll/selection state: author had scrolled down and landed on E22 --
$excel.ActiveWindow.ScrollRow = 12
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("E22").Select()

$wb.Application.Calculate()
